$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the touched cells to Text format first so numeric-looking
# strings (e.g. "7.40", "26.218.07") keep their exact characters
# instead of being auto-coerced into Number values by COM.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.218.07'
$ws.Range("E2").Value = '  -1.85%  '
$ws.Range("D3").Value = '1.583.30'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '209.71'
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("E6").Value = '  -1.26%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("D10").Value = '19.56'
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D12").Value = '1.805.75'
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("D13").Value = '1.574.19'
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").Value = '64.71'
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").Value = '26.244.42'
$ws.Range("E17").Value = '  -1.66%  '
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").Value = '7.21'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").Value = '206.69'
$ws.Range("E21").Value = '  -1.68%  '
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("E23").Value = '  -3.40%  '
$ws.Range("E24").Value = '  -1.23%  '
$ws.Range("D25").Value = '144.36'
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("D29").Value = '15.25'
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("D30").Value = '0.0503'
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("E33").Value = '  -0.73%  '
$ws.Range("D34").Value = '1.285.30'
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("E35").Value = '  +6.63%  '
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("D37").Value = '0.609'
$ws.Range("E37").Value = '  +1.10%  '
$ws.Range("E38").Value = '  -0.97%  '
$ws.Range("E39").Value = '  -1.47%  '
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").Value = '5.53'
$ws.Range("E41").Value = '  +2.27%  '
$ws.Range("D42").Value = '0.767'
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("D43").Value = '2.14'
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("D44").Value = '62.21'
$ws.Range("E44").Value = '  -1.29%  '
$ws.Range("D45").Value = '1.718.70'
$ws.Range("E45").Value = '  -1.14%  '
$ws.Range("D46").Value = '88.97'
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").Value = '7.40'
$ws.Range("E51").Value = '  +0.15%  '
